$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append new user record in row 33
$row = 33
$ws.Cells.Item($row, 1).Value = 110032
$ws.Cells.Item($row, 2).Value = 9317596770
$ws.Cells.Item($row, 3).Value = "Ewan Marsh"
$ws.Cells.Item($row, 4).Value = "ewan.marsh@xyz.com"
$ws.Cells.Item($row, 5).Value = 818876433
$ws.Cells.Item($row, 6).Value = "ACT"
$ws.Cells.Item($row, 7).Value = "eng"
$ws.Cells.Item($row, 8).Value = "PWD"
$ws.Cells.Item($row, 9).Value = $true
$ws.Cells.Item($row, 10).Value = "superadmin"
$ws.Cells.Item($row, 11).Value = "now()"
$ws.Cells.Item($row, 12).Value = "now()"

# Update selection / view to match target (clicking column M header selects the whole column)
$ws.Range("M1:XFD1048576").Select()
